# Auto-generated edit script: update column F (formation energy) values
# per commit: 'updated data from v1 of c implementation (NOTE: need to check CE calcs)'
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F20").Value = -3.124879391046773
$ws.Range("F23").Value = -3.175912136205922
$ws.Range("F24").Value = -3.192923051258973
$ws.Range("F35").Value = -3.226291917543949
$ws.Range("F36").Value = -3.223092232937459
$ws.Range("F38").Value = -3.209786069118608
$ws.Range("F39").Value = -3.208313083163581
$ws.Range("F40").Value = -3.205113398557088
$ws.Range("F41").Value = -3.200187015299124
$ws.Range("F42").Value = -3.193533933389698
$ws.Range("F43").Value = -3.186880851480266
$ws.Range("F44").Value = -3.180103329017512
$ws.Range("F45").Value = -3.173450247108085
$ws.Range("F46").Value = -3.165126226527425
$ws.Range("F47").Value = -3.163721921145488
$ws.Range("F48").Value = -3.144900347509837
$ws.Range("F49").Value = -3.14355180210814
$ws.Range("F50").Value = -3.129979004999986
$ws.Range("F51").Value = -3.123325923090555
$ws.Range("F52").Value = -3.115001902509896
$ws.Range("F53").Value = -3.110019759271698
$ws.Range("F54").Value = -3.094776023492309
$ws.Range("F55").Value = -3.087964829600126
$ws.Range("F56").Value = -3.061331743350184
$ws.Range("F57").Value = -3.051494896081042
$ws.Range("F58").Value = -3.036251160301652
$ws.Range("F78").Value = -3.308026829355379
$ws.Range("F80").Value = -3.318802654773416
$ws.Range("F81").Value = -3.326611916697479
$ws.Range("F82").Value = -3.324231653897186
$ws.Range("F83").Value = -3.329462211805341
$ws.Range("F84").Value = -3.330960694550858
$ws.Range("F85").Value = -3.332087193054831
$ws.Range("F86").Value = -3.348505084917862
$ws.Range("F87").Value = -3.336170274395798
$ws.Range("F88").Value = -3.345969226034216
$ws.Range("F89").Value = -3.352753240094705
$ws.Range("F90").Value = -3.34496102940502
$ws.Range("F93").Value = -3.352453227406527
$ws.Range("F94").Value = -3.360302306935507
$ws.Range("F95").Value = -3.362621175593469
$ws.Range("F97").Value = -3.360021413891155
$ws.Range("F101").Value = -3.369883543405971
$ws.Range("F157").Value = -3.31690202417031
$ws.Range("F159").Value = -3.307646025517744
$ws.Range("F161").Value = -3.30087809371065
$ws.Range("F162").Value = -3.298411755203607
$ws.Range("F163").Value = -3.293413858803531
$ws.Range("F164").Value = -3.289164080860678
$ws.Range("F165").Value = -3.286121245937419
$ws.Range("F166").Value = -3.278911846115435
$ws.Range("F167").Value = -3.27591042865137
$ws.Range("F168").Value = -3.275128336451912
$ws.Range("F170").Value = -3.26733727640985
$ws.Range("F171").Value = -3.261392578054374
$ws.Range("F172").Value = -3.259228449089012
$ws.Range("F173").Value = -3.255099088696018
$ws.Range("F174").Value = -3.247694263209525
$ws.Range("F175").Value = -3.247789690518993
$ws.Range("F176").Value = -3.241956785620352
$ws.Range("F177").Value = -3.236373912695765
$ws.Range("F178").Value = -3.234167422442597
$ws.Range("F179").Value = -3.226848166752846
$ws.Range("F180").Value = -3.222523481472382
$ws.Range("F181").Value = -3.216468872555873
$ws.Range("F183").Value = -3.208910814332524
$ws.Range("F184").Value = -3.200872911716354
$ws.Range("F185").Value = -3.198494678612338
$ws.Range("F186").Value = -3.194502891198111
$ws.Range("F187").Value = -3.189942772410006
$ws.Range("F188").Value = -3.183811291258322
$ws.Range("F189").Value = -3.1783107003409
$ws.Range("F191").Value = -3.168536593678417
$ws.Range("F192").Value = -3.165280881921249
$ws.Range("F193").Value = -3.159851154813457
$ws.Range("F194").Value = -3.15391621588835
$ws.Range("F196").Value = -3.140994890267179
$ws.Range("F197").Value = -3.137391156037452
$ws.Range("F198").Value = -3.132441289621505
$ws.Range("F199").Value = -3.123807322031263
$ws.Range("F200").Value = -3.118967172125306
$ws.Range("F201").Value = -3.11224512084049
$ws.Range("F203").Value = -3.101708909072259
$ws.Range("F204").Value = -3.094333157998679
$ws.Range("F221").Value = -3.441664601894059
$ws.Range("F227").Value = -3.351981791835195
$ws.Range("F228").Value = -3.28362928659911
$ws.Range("F229").Value = -3.202822879720085
$ws.Range("F233").Value = -3.529827473575495
$ws.Range("F239").Value = -3.35071770730235
$ws.Range("F240").Value = -3.269888254869889
$ws.Range("F244").Value = -3.570933898988859
$ws.Range("F251").Value = -3.315158310422232
$ws.Range("F265").Value = -3.614358253031765
$ws.Range("F276").Value = -3.634963482740108
$ws.Range("F287").Value = -3.653218696480204
$ws.Range("F298").Value = -3.667528549335481
$ws.Range("F309").Value = -3.678494475042098
